$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Multiline" column header (K2) and sample multi-line value (K3)
$ws.Range("K2").Value = "Multiline"
$ws.Range("K3").Value = "Apple" + [char]10 + "Orange"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").WrapText = $true

# Fill in the previously empty Text/General demo cells
$ws.Range("D3").Value = "Hola"
$ws.Range("E3").Value = "Hallo"

$ws.Range("E3").Select()
